$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1822.2858
$ws.Range("J112").Value = 1886.875
$ws.Range("L112").Value = 5660.625
$ws.Range("N112").Value = -7876.625
$ws.Range("H132").Value = 32803.438
$ws.Range("I132").Value = 1657
$ws.Range("J132").Value = 500000
$ws.Range("K132").Value = 4971
$ws.Range("L132").Value = 1500000
$ws.Range("M132").Value = -2441
$ws.Range("N132").Value = -1505060
$ws.Range("H135").Value = 20003314
$ws.Range("I135").Value = 25004068
$ws.Range("J135").Value = 300
$ws.Range("K135").Value = 225036612
$ws.Range("L135").Value = 2700
$ws.Range("M135").Value = -225034077
$ws.Range("N135").Value = -7770
$ws.Range("H137").Value = 3730.7144
$ws.Range("I137").Value = 3252.353
$ws.Range("K137").Value = 9757.059000000001
$ws.Range("M137").Value = -7207.059000000001
$ws.Range("H138").Value = 5211325.5
$ws.Range("J138").Value = 7579570.5
$ws.Range("L138").Value = 22738711.5
$ws.Range("N138").Value = -22748991.5
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280
$ws.Range("H141").Value = 3511.8572
$ws.Range("I141").Value = 3397.3845
$ws.Range("K141").Value = 10192.1535
$ws.Range("M141").Value = -5012.1535
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 227.66667
$ws.Range("J4").Value = 189
$ws.Range("L4").Value = 189
$ws.Range("N4").Value = -421
$ws.Range("H32").Value = 12136.508
$ws.Range("I32").Value = 10375.893
$ws.Range("K32").Value = 10375.893
$ws.Range("M32").Value = -10088.893
$ws.Range("H61").Value = 4378.129
$ws.Range("I61").Value = 4030.1428
$ws.Range("J61").Value = 4664.706
$ws.Range("K61").Value = 4030.1428
$ws.Range("L61").Value = 4664.706
$ws.Range("M61").Value = -3818.1428
$ws.Range("N61").Value = -5088.706
$ws.Range("H63").Value = 3870.875
$ws.Range("I63").Value = 3995.2856
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 3995.2856
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -3309.2856
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 3870.875
$ws.Range("I66").Value = 3995.2856
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 19976.428
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -16544.428
$ws.Range("N66").Value = -21864
$ws.Range("H74").Value = 2330.2964
$ws.Range("I74").Value = 2330.2964
$ws.Range("K74").Value = 2330.2964
$ws.Range("M74").Value = -1456.2964
$ws.Range("H77").Value = 2330.2964
$ws.Range("I77").Value = 2330.2964
$ws.Range("K77").Value = 11651.482
$ws.Range("M77").Value = -7283.482
$ws.Range("H110").Value = 3979.2856
$ws.Range("I110").Value = 3979.2856
$ws.Range("K110").Value = 3979.2856
$ws.Range("M110").Value = -1934.2856
$ws.Range("H122").Value = 3388.8823
$ws.Range("I122").Value = 1881.5349
$ws.Range("K122").Value = 5644.6047
$ws.Range("M122").Value = -3194.6047
$ws.Range("H132").Value = 4064.4255
$ws.Range("I132").Value = 3715.7
$ws.Range("K132").Value = 11147.1
$ws.Range("M132").Value = -8617.099999999999
$ws.Range("H136").Value = 4378.129
$ws.Range("I136").Value = 4030.1428
$ws.Range("J136").Value = 4664.706
$ws.Range("K136").Value = 12090.4284
$ws.Range("L136").Value = 13994.118
$ws.Range("M136").Value = -9540.428400000001
$ws.Range("N136").Value = -19094.118
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4452.3
$ws.Range("I94").Value = 4379.4
$ws.Range("J94").Value = 4525.2
$ws.Range("K94").Value = 4379.4
$ws.Range("L94").Value = 4525.2
$ws.Range("M94").Value = -3928.4
$ws.Range("N94").Value = -5427.2
$ws.Range("H105").Value = 2311.4243
$ws.Range("I105").Value = 1751.4
$ws.Range("K105").Value = 1751.4
$ws.Range("M105").Value = -4.400000000000091
$ws.Range("H119").Value = 35253.332
$ws.Range("J119").Value = 35253.332
$ws.Range("L119").Value = 35253.332
$ws.Range("N119").Value = -44929.332
$ws.Range("H126").Value = 49999.715
$ws.Range("J126").Value = 49999.715
$ws.Range("L126").Value = 49999.715
$ws.Range("N126").Value = -59879.715
$ws.Range("H134").Value = 3090.907
$ws.Range("I134").Value = 2079.2646
$ws.Range("K134").Value = 6237.793799999999
$ws.Range("M134").Value = -3702.793799999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 30799.4
$ws.Range("J28").Value = 30799.4
$ws.Range("L28").Value = 30799.4
$ws.Range("N28").Value = -31289.4
$ws.Range("H58").Value = 1749.5
$ws.Range("I58").Value = 1500
$ws.Range("K58").Value = 1500
$ws.Range("M58").Value = -1297
$ws.Range("H134").Value = 4563.1816
$ws.Range("I134").Value = 3359.5334
$ws.Range("K134").Value = 10078.6002
$ws.Range("M134").Value = -7543.600199999999
$ws.Range("H136").Value = 1749.5
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950
$ws.Range("H137").Value = 63112
$ws.Range("I137").Value = 44999.5
$ws.Range("J137").Value = 69149.5
$ws.Range("K137").Value = 44999.5
$ws.Range("L137").Value = 69149.5
$ws.Range("M137").Value = -39899.5
$ws.Range("N137").Value = -79349.5
$ws.Range("H141").Value = 296999.9
$ws.Range("J141").Value = 296999.9
$ws.Range("L141").Value = 296999.9
$ws.Range("N141").Value = -307359.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7999.8184
$ws.Range("I56").Value = 7999.8184
$ws.Range("K56").Value = 7999.8184
$ws.Range("M56").Value = -7469.8184
$ws.Range("H110").Value = 8199.799999999999
$ws.Range("I110").Value = 6999.5
$ws.Range("K110").Value = 20998.5
$ws.Range("M110").Value = -16908.5
$ws.Range("H113").Value = 1138.4
$ws.Range("J113").Value = 1056.1333
$ws.Range("L113").Value = 3168.3999
$ws.Range("N113").Value = -7508.3999
$ws.Range("H122").Value = 1703.409
$ws.Range("I122").Value = 1351.5
$ws.Range("K122").Value = 12163.5
$ws.Range("M122").Value = -9713.5
$ws.Range("H131").Value = 4308.048
$ws.Range("I131").Value = 5083.8
$ws.Range("K131").Value = 15251.4
$ws.Range("M131").Value = -10211.4
$ws.Range("H137").Value = 1469.4286
$ws.Range("I137").Value = 1469.4286
$ws.Range("K137").Value = 4408.2858
$ws.Range("M137").Value = 691.7142000000003
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 36655.305
$ws.Range("I80").Value = 64561.332
$ws.Range("J80").Value = 3168.0667
$ws.Range("K80").Value = 64561.332
$ws.Range("L80").Value = 3168.0667
$ws.Range("M80").Value = -63563.332
$ws.Range("N80").Value = -5164.066699999999
$ws.Range("H83").Value = 36655.305
$ws.Range("I83").Value = 64561.332
$ws.Range("J83").Value = 3168.0667
$ws.Range("K83").Value = 322806.66
$ws.Range("L83").Value = 15840.3335
$ws.Range("M83").Value = -317814.66
$ws.Range("N83").Value = -25824.3335
$ws.Range("H122").Value = 3072.1667
$ws.Range("I122").Value = 2677.077
$ws.Range("K122").Value = 8031.231000000001
$ws.Range("M122").Value = -5581.231000000001
$ws.Range("H132").Value = 4205.5
$ws.Range("I132").Value = 3679.8948
$ws.Range("K132").Value = 11039.6844
$ws.Range("M132").Value = -8509.6844
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2991.2334
$ws.Range("I7").Value = 1696.6
$ws.Range("K7").Value = 1696.6
$ws.Range("M7").Value = -1584.6
$ws.Range("H100").Value = 83699.21000000001
$ws.Range("I100").Value = 279377.5
$ws.Range("J100").Value = 5427.9
$ws.Range("K100").Value = 279377.5
$ws.Range("L100").Value = 5427.9
$ws.Range("M100").Value = -278836.5
$ws.Range("N100").Value = -6509.9
$ws.Range("H122").Value = 4258.6284
$ws.Range("I122").Value = 3390.8928
$ws.Range("K122").Value = 10172.6784
$ws.Range("M122").Value = -7722.678400000001
$ws.Range("H126").Value = 2991.2334
$ws.Range("I126").Value = 1696.6
$ws.Range("K126").Value = 5089.799999999999
$ws.Range("M126").Value = -2619.799999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3032.7334
$ws.Range("I81").Value = 2754.6667
$ws.Range("J81").Value = 3449.8333
$ws.Range("K81").Value = 5509.3334
$ws.Range("L81").Value = 6899.6666
$ws.Range("M81").Value = -4448.3334
$ws.Range("N81").Value = -9021.6666
$ws.Range("H84").Value = 3032.7334
$ws.Range("I84").Value = 2754.6667
$ws.Range("J84").Value = 3449.8333
$ws.Range("K84").Value = 27546.667
$ws.Range("L84").Value = 34498.333
$ws.Range("M84").Value = -22242.667
$ws.Range("N84").Value = -45106.333
$ws.Range("H100").Value = 689.381
$ws.Range("I100").Value = 689.381
$ws.Range("K100").Value = 1378.762
$ws.Range("M100").Value = -837.7619999999999
$ws.Range("H132").Value = 3587.8948
$ws.Range("I132").Value = 3412.697
$ws.Range("K132").Value = 10238.091
$ws.Range("M132").Value = -7708.091
